$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E5").Value = "seen"
